$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the date/value content of A6:A9 and the shared-string content of B6:B9
# (the underlying number format / style on A6:A9 is preserved by ClearContents).
$ws.Range("A6:B9").ClearContents()

# Update the active selection to match the edited workbook (F9).
$ws.Range("F9").Select()
